# Apply "Fetch and display dummy data" edit to the worktime tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: fill in start/end time and description (B13, C13, E13)
$ws.Range("B13").Value = 11
$ws.Range("C13").Value = 12
$ws.Range("E13").Value = "Fixed broken build and worked on R&D course"

# Row 12: update description text (E12)
$ws.Range("E12").Value = "Looking into the Flutter http package and testing example http calls to dummy JSON data"

# Row 14: fill in start/end time and description (B14, C14, E14)
$ws.Range("B14").Value = 11
$ws.Range("C14").Value = 15
$ws.Range("E14").Value = "Figuring out how to map example JSON data. Also fixing endpoints for item + character retrieval from Blizzard API"

# Row 15: fill in start/end time (B15, C15) - description stays blank
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 16

# Add a thin border around E1 (style change observed in diff)
$ws.Range("E1").Borders.Color = 0
$ws.Range("E1").Borders.LineStyle = 1

# Update selection to G13
$ws.Range("G13").Select()
